$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("CRP")
$ws4 = $wb.Worksheets.Item("CUL")
$ws5 = $wb.Worksheets.Item("GSM")
$ws6 = $wb.Worksheets.Item("LTW")
$ws7 = $wb.Worksheets.Item("WVR")

# ALC!row9
$ws1.Range("H9").Value = 147.42857
$ws1.Range("I9").Value = 142.4
$ws1.Range("J9").Value = 160.0
$ws1.Range("K9").Value = 142.4
$ws1.Range("L9").Value = 160.0
$ws1.Range("M9").Value = 26.59999999999999
$ws1.Range("N9").Value = -498.0

# ALC!row40
$ws1.Range("H40").Value = 2942.8572
$ws1.Range("I40").Value = 3400.0
$ws1.Range("J40").Value = 1800.0
$ws1.Range("K40").Value = 3400.0
$ws1.Range("L40").Value = 1800.0
$ws1.Range("M40").Value = -3225.0
$ws1.Range("N40").Value = -2150.0

# ALC!row43
$ws1.Range("H43").Value = 6185638.5
$ws1.Range("J43").Value = 27778278.0
$ws1.Range("L43").Value = 27778278.0
$ws1.Range("N43").Value = -27778416.0

# ALC!row138
$ws1.Range("H138").Value = 436574.6
$ws1.Range("I138").Value = 1865.5
$ws1.Range("J138").Value = 519376.3
$ws1.Range("K138").Value = 5596.5
$ws1.Range("L138").Value = 1558128.9
$ws1.Range("M138").Value = -456.5
$ws1.Range("N138").Value = -1568408.9

# ARM!row88
$ws2.Range("H88").Value = 2751.4285
$ws2.Range("I88").Value = 0.0
$ws2.Range("K88").Value = 0.0
$ws2.Range("M88").ClearContents()

# ARM!row91
$ws2.Range("H91").Value = 2751.4285
$ws2.Range("I91").Value = 0.0
$ws2.Range("K91").Value = 0.0
$ws2.Range("M91").ClearContents()

# CRP!row114
$ws3.Range("H114").Value = 26660.0
$ws3.Range("I114").Value = 21000.0
$ws3.Range("K114").Value = 21000.0
$ws3.Range("M114").Value = -16661.0

# CUL!row34
$ws4.Range("H34").Value = 2588.3
$ws4.Range("I34").Value = 1723.0
$ws4.Range("J34").Value = 3165.1667
$ws4.Range("K34").Value = 5169.0
$ws4.Range("L34").Value = 9495.500100000001
$ws4.Range("M34").Value = -5085.0
$ws4.Range("N34").Value = -9663.500100000001

# CUL!row39
$ws4.Range("H39").Value = 2501.3157
$ws4.Range("J39").Value = 2257.8125
$ws4.Range("L39").Value = 6773.4375
$ws4.Range("N39").Value = -7361.4375

# CUL!row55
$ws4.Range("H55").Value = 3468.5715
$ws4.Range("I55").Value = 0.0
$ws4.Range("J55").Value = 3468.5715
$ws4.Range("K55").Value = 0.0
$ws4.Range("L55").Value = 10405.7145
$ws4.Range("M55").ClearContents()
$ws4.Range("N55").Value = -10759.7145

# CUL!row68
$ws4.Range("H68").Value = 1260.6
$ws4.Range("I68").Value = 1086.0
$ws4.Range("J68").Value = 1668.0
$ws4.Range("K68").Value = 3258.0
$ws4.Range("L68").Value = 5004.0
$ws4.Range("M68").Value = -2447.0
$ws4.Range("N68").Value = -6626.0

# CUL!row69
$ws4.Range("H69").Value = 2737.6365
$ws4.Range("I69").Value = 0.0
$ws4.Range("J69").Value = 2737.6365
$ws4.Range("K69").Value = 0.0
$ws4.Range("L69").Value = 8212.9095
$ws4.Range("M69").ClearContents()
$ws4.Range("N69").Value = -9834.9095

# CUL!row70
$ws4.Range("H70").Value = 3989.1667
$ws4.Range("I70").Value = 1470.0
$ws4.Range("J70").Value = 5248.75
$ws4.Range("K70").Value = 4410.0
$ws4.Range("L70").Value = 15746.25
$ws4.Range("M70").Value = -4095.0
$ws4.Range("N70").Value = -16376.25

# CUL!row71
$ws4.Range("H71").Value = 1260.6
$ws4.Range("I71").Value = 1086.0
$ws4.Range("J71").Value = 1668.0
$ws4.Range("K71").Value = 9774.0
$ws4.Range("L71").Value = 15012.0
$ws4.Range("M71").Value = -5718.0
$ws4.Range("N71").Value = -23124.0

# CUL!row72
$ws4.Range("H72").Value = 2737.6365
$ws4.Range("I72").Value = 0.0
$ws4.Range("J72").Value = 2737.6365
$ws4.Range("K72").Value = 0.0
$ws4.Range("L72").Value = 24638.7285
$ws4.Range("M72").ClearContents()
$ws4.Range("N72").Value = -32750.7285

# CUL!row73
$ws4.Range("H73").Value = 3989.1667
$ws4.Range("I73").Value = 1470.0
$ws4.Range("J73").Value = 5248.75
$ws4.Range("K73").Value = 4410.0
$ws4.Range("L73").Value = 15746.25
$ws4.Range("M73").Value = -3318.0
$ws4.Range("N73").Value = -17930.25

# CUL!row82
$ws4.Range("H82").Value = 6160.2
$ws4.Range("I82").Value = 2266.6667
$ws4.Range("J82").Value = 12000.5
$ws4.Range("K82").Value = 6800.000100000001
$ws4.Range("L82").Value = 36001.5
$ws4.Range("M82").Value = -6394.000100000001
$ws4.Range("N82").Value = -36813.5

# CUL!row85
$ws4.Range("H85").Value = 6160.2
$ws4.Range("I85").Value = 2266.6667
$ws4.Range("J85").Value = 12000.5
$ws4.Range("K85").Value = 6800.000100000001
$ws4.Range("L85").Value = 36001.5
$ws4.Range("M85").Value = -5396.000100000001
$ws4.Range("N85").Value = -38809.5

# CUL!row92
$ws4.Range("H92").Value = 283.05264
$ws4.Range("I92").Value = 354.66666
$ws4.Range("J92").Value = 250.0
$ws4.Range("K92").Value = 1063.99998
$ws4.Range("L92").Value = 750.0
$ws4.Range("M92").Value = 184.0000199999999
$ws4.Range("N92").Value = -3246.0

# CUL!row93
$ws4.Range("H93").Value = 6505.4
$ws4.Range("J93").Value = 6505.4
$ws4.Range("L93").Value = 19516.2
$ws4.Range("N93").Value = -23260.2

# CUL!row94
$ws4.Range("H94").Value = 4384.857
$ws4.Range("J94").Value = 4278.3335
$ws4.Range("L94").Value = 12835.0005
$ws4.Range("N94").Value = -14187.0005

# CUL!row98
$ws4.Range("H98").Value = 359.2
$ws4.Range("I98").Value = 359.2
$ws4.Range("J98").Value = 0.0
$ws4.Range("K98").Value = 1077.6
$ws4.Range("L98").Value = 0.0
$ws4.Range("M98").Value = 420.4000000000001
$ws4.Range("N98").ClearContents()

# CUL!row99
$ws4.Range("H99").Value = 1811.7778
$ws4.Range("I99").Value = 400.0
$ws4.Range("J99").Value = 2517.6667
$ws4.Range("K99").Value = 1200.0
$ws4.Range("L99").Value = 7553.000100000001
$ws4.Range("M99").Value = 1046.0
$ws4.Range("N99").Value = -12045.0001

# CUL!row100
$ws4.Range("H100").Value = 3473.5386
$ws4.Range("J100").Value = 3473.5386
$ws4.Range("L100").Value = 10420.6158
$ws4.Range("N100").Value = -12042.6158

# CUL!row104
$ws4.Range("H104").Value = 4743.3
$ws4.Range("I104").Value = 3526.0
$ws4.Range("J104").Value = 5554.8335
$ws4.Range("K104").Value = 10578.0
$ws4.Range("L104").Value = 16664.5005
$ws4.Range("M104").Value = -7957.0
$ws4.Range("N104").Value = -21906.5005

# CUL!row105
$ws4.Range("H105").Value = 989999.0
$ws4.Range("J105").Value = 989999.0
$ws4.Range("L105").Value = 2969997.0
$ws4.Range("N105").Value = -2975239.0

# CUL!row106
$ws4.Range("H106").Value = 2943.0
$ws4.Range("J106").Value = 2943.0
$ws4.Range("L106").Value = 8829.0
$ws4.Range("N106").Value = -10721.0

# CUL!row113
$ws4.Range("H113").Value = 682.02704
$ws4.Range("I113").Value = 600.0
$ws4.Range("J113").Value = 684.30554
$ws4.Range("K113").Value = 1800.0
$ws4.Range("L113").Value = 2052.91662
$ws4.Range("M113").Value = 370.0
$ws4.Range("N113").Value = -6392.91662

# GSM!row80
$ws5.Range("H80").Value = 6800.0
$ws5.Range("I80").Value = 0.0
$ws5.Range("K80").Value = 0.0
$ws5.Range("M80").ClearContents()

# GSM!row83
$ws5.Range("H83").Value = 6800.0
$ws5.Range("I83").Value = 0.0
$ws5.Range("K83").Value = 0.0
$ws5.Range("M83").ClearContents()

# LTW!row40
$ws6.Range("H40").Value = 4631.25
$ws6.Range("I40").Value = 3138.5715
$ws6.Range("J40").Value = 6721.0
$ws6.Range("K40").Value = 3138.5715
$ws6.Range("L40").Value = 6721.0
$ws6.Range("M40").Value = -3002.5715
$ws6.Range("N40").Value = -6993.0

# WVR!row126
$ws7.Range("H126").Value = 66668620.0
$ws7.Range("I126").Value = 76924870.0
$ws7.Range("J126").Value = 2999.5
$ws7.Range("K126").Value = 230774610.0
$ws7.Range("L126").Value = 8998.5
$ws7.Range("M126").Value = -230772140.0
$ws7.Range("N126").Value = -13938.5
